# Generate Report for Handoff
#
# The localization-status workbook tracks handoff progress for a set of
# source files across target languages (zh-cn, de-de) plus an "Overview"
# rollup sheet. This run represents a fresh handoff-report generation:
#   - the four files that were "Ready for handoff" (previously queued with
#     Priority "low") are now flagged with Priority "ht" (hand-translated)
#   - their "Latest Handoff" timestamps advance to the new generation run

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Rows 4-7 on every sheet correspond to the same four files:
#   0a7b0ba3-..., 320f2bc1-..., 551c3170-..., ec341ab7-...
$rows = 4, 5, 6, 7

foreach ($r in $rows) {
    # Overview: "Latest HO Xliff Generate Date" (column G) moves forward
    $overview.Range("G$r").Value = "2016-08-17 14:30:43"

    # zh-cn: Priority (E) low -> ht, Latest Handoff Datetime (H) advances
    $zhcn.Range("E$r").Value = "ht"
    $zhcn.Range("H$r").Value = "2016-08-17 14:30:35"

    # de-de: Priority (E) low -> ht, Latest Handoff Datetime (H) advances
    $dede.Range("E$r").Value = "ht"
    $dede.Range("H$r").Value = "2016-08-17 14:30:43"
}
